$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.41
$summary.Range("B4").Value = 1.2
$summary.Range("B6").Value = 125
$summary.Range("B7").Value = 55
$summary.Range("B9").Value = 44

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.41
$status.Range("D5").Value = 92
$status.Range("E5").Value = 1.09
$status.Range("F5").Value = 1.41
$status.Range("G5").Value = 44.57

# ---------------------------------------------------------------------------
# All Trades sheet - trade #125 (row 126) closes, and a new trade #158
# (row 159) is appended
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(126, 7).Value = 0.9
$allTrades.Cells.Item(126, 8).Value = "CLOSED"
$allTrades.Cells.Item(126, 9).Value = 4.6512
$allTrades.Cells.Item(126, 10).Value = 0.04
$allTrades.Cells.Item(126, 11).Value = 101.41
$allTrades.Cells.Item(126, 12).Value = "early_exit"
$allTrades.Cells.Item(126, 13).Value = 0.11

$allTrades.Cells.Item(159, 1).Value = 158
$allTrades.Cells.Item(159, 2).NumberFormat = "@"
$allTrades.Cells.Item(159, 2).Value = "2026-02-17"
$allTrades.Cells.Item(159, 3).Value = "21:27:03"
$allTrades.Cells.Item(159, 4).Value = "MarketMaking"
$allTrades.Cells.Item(159, 5).Value = "DOWN"
$allTrades.Cells.Item(159, 6).Value = 0.86
$allTrades.Cells.Item(159, 8).Value = "OPEN"
$allTrades.Cells.Item(159, 9).Value = 0
$allTrades.Cells.Item(159, 10).Value = 0
$allTrades.Cells.Item(159, 11).Value = 101.371797784678
$allTrades.Cells.Item(159, 13).Value = 0
$allTrades.Cells.Item(159, 14).Value = 0
$allTrades.Cells.Item(159, 15).Value = 0
$allTrades.Cells.Item(159, 16).Value = 0.6
$allTrades.Cells.Item(159, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet - trade #125 (row 93) closes, and trade #158 (row 126)
# is appended (same trades as above, mirrored into the strategy-specific tab)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Cells.Item(93, 7).Value = 0.9
$mm.Cells.Item(93, 8).Value = "CLOSED"
$mm.Cells.Item(93, 9).Value = 4.6512
$mm.Cells.Item(93, 10).Value = 0.04
$mm.Cells.Item(93, 11).Value = 101.41
$mm.Cells.Item(93, 16).Value = "early_exit"
$mm.Cells.Item(93, 17).Value = 0.11

$mm.Cells.Item(126, 1).Value = 158
$mm.Cells.Item(126, 2).NumberFormat = "@"
$mm.Cells.Item(126, 2).Value = "2026-02-17"
$mm.Cells.Item(126, 3).Value = "21:27:03"
$mm.Cells.Item(126, 4).Value = "MarketMaking"
$mm.Cells.Item(126, 5).Value = "DOWN"
$mm.Cells.Item(126, 6).Value = 0.86
$mm.Cells.Item(126, 8).Value = "OPEN"
$mm.Cells.Item(126, 9).Value = 0
$mm.Cells.Item(126, 10).Value = 0
$mm.Cells.Item(126, 11).Value = 101.371797784678
$mm.Cells.Item(126, 12).Value = 0
$mm.Cells.Item(126, 13).Value = 0
$mm.Cells.Item(126, 14).Value = 0.6
$mm.Cells.Item(126, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(126, 17).Value = 0
